$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Alunos")
$ws.Activate()

# Clear the contents (not formatting) of the student rows, keeping the
# number-format styling on column D intact.
$ws.Range("A2:F5").ClearContents()

# Leave the selection on E4, matching the saved selection in the workbook.
$ws.Range("E4").Select()
